$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 150
$ws.Range("I12").Value = 122.5
$ws.Range("J12").Value = 163.75
$ws.Range("K12").Value = 122.5
$ws.Range("L12").Value = 163.75
$ws.Range("M12").Value = 47.5
$ws.Range("N12").Value = -503.75
$ws.Range("H17").Value = 1194.5161
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 1194.5161
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 3583.5483
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -3919.5483
$ws.Range("H80").Value = 629.25
$ws.Range("I80").Value = 424.875
$ws.Range("J80").Value = 1038
$ws.Range("K80").Value = 1274.625
$ws.Range("L80").Value = 3114
$ws.Range("M80").Value = -276.625
$ws.Range("N80").Value = -5110
$ws.Range("H82").Value = 1491.7
$ws.Range("I82").Value = 1491.7
$ws.Range("K82").Value = 4475.1
$ws.Range("M82").Value = -4069.1
$ws.Range("H83").Value = 629.25
$ws.Range("I83").Value = 424.875
$ws.Range("J83").Value = 1038
$ws.Range("K83").Value = 3823.875
$ws.Range("L83").Value = 9342
$ws.Range("M83").Value = 1168.125
$ws.Range("N83").Value = -19326
$ws.Range("H85").Value = 1491.7
$ws.Range("I85").Value = 1491.7
$ws.Range("K85").Value = 4475.1
$ws.Range("M85").Value = -3071.1
$ws.Range("H87").Value = 76800.914
$ws.Range("J87").Value = 76800.914
$ws.Range("L87").Value = 76800.914
$ws.Range("N87").Value = -79296.914
$ws.Range("H88").Value = 3215.6155
$ws.Range("I88").Value = 474.5
$ws.Range("J88").Value = 3444.0417
$ws.Range("K88").Value = 474.5
$ws.Range("L88").Value = 3444.0417
$ws.Range("M88").Value = -68.5
$ws.Range("N88").Value = -4256.0417
$ws.Range("H90").Value = 76800.914
$ws.Range("J90").Value = 76800.914
$ws.Range("L90").Value = 230402.742
$ws.Range("N90").Value = -242882.742
$ws.Range("H91").Value = 3215.6155
$ws.Range("I91").Value = 474.5
$ws.Range("J91").Value = 3444.0417
$ws.Range("K91").Value = 474.5
$ws.Range("L91").Value = 3444.0417
$ws.Range("M91").Value = 929.5
$ws.Range("N91").Value = -6252.0417
$ws.Range("H106").Value = 34098
$ws.Range("I106").Value = 34098
$ws.Range("K106").Value = 34098
$ws.Range("M106").Value = -33467
$ws.Range("H112").Value = 7880.9443
$ws.Range("J112").Value = 8168.9414
$ws.Range("L112").Value = 24506.8242
$ws.Range("N112").Value = -26722.8242
$ws.Range("H116").Value = 4506
$ws.Range("I116").Value = 3525.5
$ws.Range("J116").Value = 5813.3335
$ws.Range("K116").Value = 3525.5
$ws.Range("L116").Value = 5813.3335
$ws.Range("M116").Value = -83.5
$ws.Range("N116").Value = -12697.3335
$ws.Range("H132").Value = 20836004
$ws.Range("I132").Value = 22224966
$ws.Range("K132").Value = 66674898
$ws.Range("M132").Value = -66672368
$ws.Range("H137").Value = 73247.03999999999
$ws.Range("I137").Value = 82780.95
$ws.Range("K137").Value = 248342.85
$ws.Range("M137").Value = -245792.85
$ws.Range("H138").Value = 3417.8809
$ws.Range("I138").Value = 1432.9445
$ws.Range("J138").Value = 4906.5835
$ws.Range("K138").Value = 4298.833500000001
$ws.Range("L138").Value = 14719.7505
$ws.Range("M138").Value = 841.1664999999994
$ws.Range("N138").Value = -24999.7505

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2925045.5
$ws.Range("I2").Value = 4273980
$ws.Range("K2").Value = 4273980
$ws.Range("M2").Value = -4273867
$ws.Range("H32").Value = 4813.8643
$ws.Range("I32").Value = 3331.0637
$ws.Range("K32").Value = 3331.0637
$ws.Range("M32").Value = -3044.0637
$ws.Range("H41").Value = 2922.45
$ws.Range("I41").Value = 1560.5385
$ws.Range("K41").Value = 1560.5385
$ws.Range("M41").Value = -1146.5385
$ws.Range("H61").Value = 6605.125
$ws.Range("I61").Value = 7319.15
$ws.Range("J61").Value = 3035
$ws.Range("K61").Value = 7319.15
$ws.Range("L61").Value = 3035
$ws.Range("M61").Value = -7107.15
$ws.Range("N61").Value = -3459
$ws.Range("H74").Value = 47797.562
$ws.Range("I74").Value = 10775.053
$ws.Range("J74").Value = 188483.1
$ws.Range("K74").Value = 10775.053
$ws.Range("L74").Value = 188483.1
$ws.Range("M74").Value = -9901.053
$ws.Range("N74").Value = -190231.1
$ws.Range("H77").Value = 47797.562
$ws.Range("I77").Value = 10775.053
$ws.Range("J77").Value = 188483.1
$ws.Range("K77").Value = 53875.265
$ws.Range("L77").Value = 942415.5
$ws.Range("M77").Value = -49507.265
$ws.Range("N77").Value = -951151.5
$ws.Range("H92").Value = 74772
$ws.Range("J92").Value = 74772
$ws.Range("L92").Value = 74772
$ws.Range("N92").Value = -79764
$ws.Range("H110").Value = 1544521.2
$ws.Range("I110").Value = 1737454.8
$ws.Range("K110").Value = 1737454.8
$ws.Range("M110").Value = -1735409.8
$ws.Range("H116").Value = 2925045.5
$ws.Range("I116").Value = 4273980
$ws.Range("K116").Value = 4273980
$ws.Range("M116").Value = -4271686
$ws.Range("H132").Value = 3312.8645
$ws.Range("I132").Value = 3567.4
$ws.Range("J132").Value = 2777
$ws.Range("K132").Value = 10702.2
$ws.Range("L132").Value = 8331
$ws.Range("M132").Value = -8172.200000000001
$ws.Range("N132").Value = -13391
$ws.Range("H136").Value = 6605.125
$ws.Range("I136").Value = 7319.15
$ws.Range("J136").Value = 3035
$ws.Range("K136").Value = 21957.45
$ws.Range("L136").Value = 9105
$ws.Range("M136").Value = -19407.45
$ws.Range("N136").Value = -14205

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2925045.5
$ws.Range("I3").Value = 4273980
$ws.Range("K3").Value = 4273980
$ws.Range("M3").Value = -4273866
$ws.Range("H8").Value = 476.75
$ws.Range("I8").Value = 476.75
$ws.Range("K8").Value = 476.75
$ws.Range("M8").Value = -336.75
$ws.Range("H94").Value = 9099698
$ws.Range("I94").Value = 18182218
$ws.Range("K94").Value = 18182218
$ws.Range("M94").Value = -18181767
$ws.Range("H99").Value = 4204068
$ws.Range("I99").Value = 5496654
$ws.Range("K99").Value = 5496654
$ws.Range("M99").Value = -5495156
$ws.Range("H134").Value = 10741.061
$ws.Range("I134").Value = 10722.692
$ws.Range("J134").Value = 10809.286
$ws.Range("K134").Value = 32168.076
$ws.Range("L134").Value = 32427.858
$ws.Range("M134").Value = -29633.076
$ws.Range("N134").Value = -37497.858

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4096.6
$ws.Range("I31").Value = 1586.2593
$ws.Range("J31").Value = 5025.082
$ws.Range("K31").Value = 1586.2593
$ws.Range("L31").Value = 5025.082
$ws.Range("M31").Value = -1291.2593
$ws.Range("N31").Value = -5615.082
$ws.Range("H34").Value = 4096.6
$ws.Range("I34").Value = 1586.2593
$ws.Range("J34").Value = 5025.082
$ws.Range("K34").Value = 1586.2593
$ws.Range("L34").Value = 5025.082
$ws.Range("M34").Value = -1384.2593
$ws.Range("N34").Value = -5429.082
$ws.Range("H58").Value = 3232.4707
$ws.Range("I58").Value = 2873.5833
$ws.Range("J58").Value = 4093.8
$ws.Range("K58").Value = 2873.5833
$ws.Range("L58").Value = 4093.8
$ws.Range("M58").Value = -2670.5833
$ws.Range("N58").Value = -4499.8
$ws.Range("H74").Value = 51662.668
$ws.Range("J74").Value = 51662.668
$ws.Range("L74").Value = 51662.668
$ws.Range("N74").Value = -53410.668
$ws.Range("H77").Value = 51662.668
$ws.Range("J77").Value = 51662.668
$ws.Range("L77").Value = 154988.004
$ws.Range("N77").Value = -163724.004
$ws.Range("H94").Value = 757
$ws.Range("J94").Value = 891.58826
$ws.Range("L94").Value = 891.58826
$ws.Range("N94").Value = -1793.58826
$ws.Range("H105").Value = 1853.4546
$ws.Range("I105").Value = 1486
$ws.Range("K105").Value = 1486
$ws.Range("M105").Value = 261
$ws.Range("H132").Value = 74280.21000000001
$ws.Range("I132").Value = 79686.38
$ws.Range("K132").Value = 239059.14
$ws.Range("M132").Value = -236529.14
$ws.Range("H134").Value = 12011.261
$ws.Range("I134").Value = 11079.917
$ws.Range("J134").Value = 13027.272
$ws.Range("K134").Value = 33239.751
$ws.Range("L134").Value = 39081.81600000001
$ws.Range("M134").Value = -30704.751
$ws.Range("N134").Value = -44151.81600000001
$ws.Range("H136").Value = 3232.4707
$ws.Range("I136").Value = 2873.5833
$ws.Range("J136").Value = 4093.8
$ws.Range("K136").Value = 8620.749899999999
$ws.Range("L136").Value = 12281.4
$ws.Range("M136").Value = -6070.749899999999
$ws.Range("N136").Value = -17381.4

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 734.3333
$ws.Range("I8").Value = 734.3333
$ws.Range("K8").Value = 2202.9999
$ws.Range("M8").Value = -2063.9999
$ws.Range("H34").Value = 414.5
$ws.Range("I34").Value = 97.40000000000001
$ws.Range("J34").Value = 2000
$ws.Range("K34").Value = 292.2
$ws.Range("L34").Value = 6000
$ws.Range("M34").Value = -208.2
$ws.Range("N34").Value = -6168
$ws.Range("H39").Value = 2233.6667
$ws.Range("J39").Value = 1750.75
$ws.Range("L39").Value = 5252.25
$ws.Range("N39").Value = -5840.25
$ws.Range("H55").Value = 47180.047
$ws.Range("J55").Value = 128685
$ws.Range("L55").Value = 386055
$ws.Range("N55").Value = -386409
$ws.Range("H69").Value = 4649.6
$ws.Range("I69").Value = 4583
$ws.Range("K69").Value = 13749
$ws.Range("M69").Value = -12938
$ws.Range("H72").Value = 4649.6
$ws.Range("I72").Value = 4583
$ws.Range("K72").Value = 41247
$ws.Range("M72").Value = -37191
$ws.Range("H93").Value = 5000
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H109").Value = 6845
$ws.Range("I109").Value = 7589.636
$ws.Range("J109").Value = 2749.5
$ws.Range("K109").Value = 22768.908
$ws.Range("L109").Value = 8248.5
$ws.Range("M109").Value = -21728.908
$ws.Range("N109").Value = -10328.5
$ws.Range("H121").Value = 270.33334
$ws.Range("I121").Value = 269.42856
$ws.Range("J121").Value = 273.5
$ws.Range("K121").Value = 808.28568
$ws.Range("L121").Value = 820.5
$ws.Range("M121").Value = 501.71432
$ws.Range("N121").Value = -3440.5
$ws.Range("H128").Value = 179374.75
$ws.Range("I128").Value = 179374.75
$ws.Range("K128").Value = 538124.25
$ws.Range("M128").Value = -533144.25
$ws.Range("H131").Value = 24512652
$ws.Range("J131").Value = 20836358
$ws.Range("L131").Value = 62509074
$ws.Range("N131").Value = -62519154
$ws.Range("H132").Value = 920
$ws.Range("I132").Value = 900
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 8100
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -5570
$ws.Range("N132").Value = -14060
$ws.Range("H136").Value = 2600
$ws.Range("H140").Value = 2774.875
$ws.Range("I140").Value = 2774.875
$ws.Range("K140").Value = 8324.625
$ws.Range("M140").Value = -3144.625

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 10110.263
$ws.Range("I19").Value = 10005
$ws.Range("J19").Value = 10405
$ws.Range("K19").Value = 10005
$ws.Range("L19").Value = 10405
$ws.Range("M19").Value = -9717
$ws.Range("N19").Value = -10981
$ws.Range("H70").Value = 10531263
$ws.Range("I70").Value = 15388976
$ws.Range("K70").Value = 15388976
$ws.Range("M70").Value = -15388706
$ws.Range("H73").Value = 10531263
$ws.Range("I73").Value = 15388976
$ws.Range("K73").Value = 15388976
$ws.Range("M73").Value = -15388040
$ws.Range("H80").Value = 1752393.1
$ws.Range("I80").Value = 2450764.8
$ws.Range("J80").Value = 6464
$ws.Range("K80").Value = 2450764.8
$ws.Range("L80").Value = 6464
$ws.Range("M80").Value = -2449766.8
$ws.Range("N80").Value = -8460
$ws.Range("H83").Value = 1752393.1
$ws.Range("I83").Value = 2450764.8
$ws.Range("J83").Value = 6464
$ws.Range("K83").Value = 12253824
$ws.Range("L83").Value = 32320
$ws.Range("M83").Value = -12248832
$ws.Range("N83").Value = -42304
$ws.Range("H113").Value = 10421117
$ws.Range("I113").Value = 23813270
$ws.Range("K113").Value = 23813270
$ws.Range("M113").Value = -23811100
$ws.Range("H122").Value = 1112447.2
$ws.Range("I122").Value = 1482763
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 4448289
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -4445839
$ws.Range("N122").Value = -9400
$ws.Range("H132").Value = 6115.3774
$ws.Range("I132").Value = 4273.1333
$ws.Range("J132").Value = 16478
$ws.Range("K132").Value = 12819.3999
$ws.Range("L132").Value = 49434
$ws.Range("M132").Value = -10289.3999
$ws.Range("N132").Value = -54494

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H41").Value = 15000
$ws.Range("J41").Value = 15000
$ws.Range("L41").Value = 15000
$ws.Range("N41").Value = -15876
$ws.Range("H45").Value = 15735.625
$ws.Range("J45").Value = 9523
$ws.Range("L45").Value = 9523
$ws.Range("N45").Value = -10337
$ws.Range("H59").Value = 35299.332
$ws.Range("J59").Value = 35299.332
$ws.Range("L59").Value = 35299.332
$ws.Range("N59").Value = -36607.332
$ws.Range("H100").Value = 3659.5925
$ws.Range("I100").Value = 3608.2693
$ws.Range("K100").Value = 3608.2693
$ws.Range("M100").Value = -3067.2693
$ws.Range("H104").Value = 9479.200000000001
$ws.Range("J104").Value = 9479.200000000001
$ws.Range("L104").Value = 9479.200000000001
$ws.Range("N104").Value = -16467.2
$ws.Range("H132").Value = 11297.772
$ws.Range("I132").Value = 11620.706
$ws.Range("K132").Value = 34862.118
$ws.Range("M132").Value = -32332.118
$ws.Range("H136").Value = 45412.36
$ws.Range("I136").Value = 54765.5
$ws.Range("K136").Value = 164296.5
$ws.Range("M136").Value = -161746.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 49995
$ws.Range("J2").Value = 49995
$ws.Range("L2").Value = 49995
$ws.Range("N2").Value = -50219
$ws.Range("H4").Value = 9994.625
$ws.Range("I4").Value = 9994
$ws.Range("J4").Value = 9996.5
$ws.Range("K4").Value = 9994
$ws.Range("L4").Value = 9996.5
$ws.Range("M4").Value = -9881
$ws.Range("N4").Value = -10222.5
$ws.Range("H14").Value = 4000
$ws.Range("I14").Value = 3000
$ws.Range("J14").Value = 5000
$ws.Range("K14").Value = 3000
$ws.Range("L14").Value = 5000
$ws.Range("M14").Value = -2832
$ws.Range("N14").Value = -5336
$ws.Range("H48").Value = 30376.666
$ws.Range("I48").Value = 25000
$ws.Range("K48").Value = 25000
$ws.Range("M48").Value = -24431
$ws.Range("H107").Value = 45456596
$ws.Range("I107").Value = 58825996
$ws.Range("K107").Value = 176477988
$ws.Range("M107").Value = -176476068
$ws.Range("H122").Value = 2940.72
$ws.Range("I122").Value = 1981.7142
$ws.Range("J122").Value = 4161.273
$ws.Range("K122").Value = 5945.142599999999
$ws.Range("L122").Value = 12483.819
$ws.Range("M122").Value = -3495.142599999999
$ws.Range("N122").Value = -17383.819
$ws.Range("H132").Value = 23517898
$ws.Range("I132").Value = 30311070
$ws.Range("K132").Value = 90933210
$ws.Range("M132").Value = -90930680
$ws.Range("H136").Value = 9461.809999999999
$ws.Range("I136").Value = 10188.177
$ws.Range("K136").Value = 30564.531
$ws.Range("M136").Value = -28014.531
